$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at position 17 (shifts rows 17-28 down to 18-29)
$ws1.Rows.Item(17).Insert()

# New row 17 content (brand-new con event inserted)
$ws1.Cells.Item(17, 3).Value = '上海·第五十四届妖漫动漫游戏展'
$ws1.Cells.Item(17, 4).Value = '秀浦路668号 新田360广场(上海康桥店)'
$ws1.Cells.Item(17, 6).Value = 0
$ws1.Cells.Item(17, 7).Value = 80
$ws1.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82097'
$ws1.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/GNchtawR1708938307099.jpeg'
$ws1.Cells.Item(17, 1).Value = 16

# Update "want-to-go" counts (F column) for rows shifted down by the insert
$ws1.Cells.Item(18, 6).Value = 1392
$ws1.Cells.Item(19, 6).Value = 624
$ws1.Cells.Item(20, 6).Value = 363
$ws1.Cells.Item(21, 6).Value = 81
$ws1.Cells.Item(22, 6).Value = 1041
$ws1.Cells.Item(23, 6).Value = 92
$ws1.Cells.Item(24, 6).Value = 2111
$ws1.Cells.Item(25, 6).Value = 213
$ws1.Cells.Item(26, 6).Value = 61
$ws1.Cells.Item(27, 6).Value = 374
$ws1.Cells.Item(28, 6).Value = 51
$ws1.Cells.Item(29, 6).Value = 3428

# Other in-place cell updates on sheet 1 (rows not affected by the shift)
$ws1.Cells.Item(2, 6).Value = 12
$ws1.Cells.Item(5, 6).Value = 6169
$ws1.Cells.Item(6, 6).Value = 687
$ws1.Cells.Item(7, 6).Value = 1072
$ws1.Cells.Item(8, 6).Value = 43
$ws1.Cells.Item(10, 6).Value = 288
$ws1.Cells.Item(12, 6).Value = 629
$ws1.Cells.Item(13, 3).Value = '上海·城市动漫节'
$ws1.Cells.Item(13, 6).Value = 1076
$ws1.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/OahoM3s41708933408573.jpeg'
$ws1.Cells.Item(14, 6).Value = 66
$ws1.Cells.Item(16, 6).Value = 339

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(5, 6).Value = 116
$ws2.Cells.Item(9, 6).Value = 675
$ws2.Cells.Item(14, 6).Value = 93
$ws2.Cells.Item(15, 6).Value = 633
$ws2.Cells.Item(20, 6).Value = 307
$ws2.Cells.Item(21, 6).Value = 4073
$ws2.Cells.Item(25, 6).Value = 167
$ws2.Cells.Item(29, 6).Value = 202
$ws2.Cells.Item(33, 7).Value = '不可售'

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(9, 6).Value = 424
$ws3.Cells.Item(10, 6).Value = 125
$ws3.Cells.Item(12, 6).Value = 705

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(7, 6).Value = 424
$ws4.Cells.Item(8, 6).Value = 125
$ws4.Cells.Item(9, 6).Value = 705
$ws4.Cells.Item(10, 6).Value = 116
$ws4.Cells.Item(13, 6).Value = 6169
$ws4.Cells.Item(15, 6).Value = 687
$ws4.Cells.Item(16, 6).Value = 1072
$ws4.Cells.Item(17, 6).Value = 675
$ws4.Cells.Item(18, 6).Value = 43
$ws4.Cells.Item(20, 6).Value = 288
$ws4.Cells.Item(22, 6).Value = 629
$ws4.Cells.Item(24, 6).Value = 93
$ws4.Cells.Item(27, 3).Value = '上海·城市动漫节'
$ws4.Cells.Item(27, 6).Value = 1076
$ws4.Cells.Item(27, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/OahoM3s41708933408573.jpeg'
$ws4.Cells.Item(28, 6).Value = 66
$ws4.Cells.Item(30, 6).Value = 339
$ws4.Cells.Item(33, 6).Value = 1392
$ws4.Cells.Item(34, 6).Value = 624
$ws4.Cells.Item(35, 6).Value = 363
$ws4.Cells.Item(37, 6).Value = 167
$ws4.Cells.Item(39, 6).Value = 1041
$ws4.Cells.Item(40, 6).Value = 92
$ws4.Cells.Item(42, 6).Value = 2111
$ws4.Cells.Item(45, 6).Value = 61
$ws4.Cells.Item(47, 6).Value = 51
$ws4.Cells.Item(48, 6).Value = 3428

